$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Text change: "Projektin käynnistys" (E4) becomes "Frontin käynnistys" ---
$ws.Range("E4").Value = "Frontin käynnistys"

# --- Capture the ORIGINAL E8 formatting (plain/default style) for reuse on
#     the new E11 cell before E8 itself gets restyled below. ---
$ws.Range("E8").Copy()
$ws.Range("E11").PasteSpecial($fmt)

# --- New row 10: Bäkkärin käynnistys (values/formulas first, then style) ---
$ws.Range("A10").Value = 46034
$ws.Range("B10").Value = 0.65625
$ws.Range("C10").Value = 0.6875
$ws.Range("D10").Formula = "=C10-B10"
$ws.Range("E10").Value = "Bäkkärin käynnistys"

# --- New row 11: Bäkkäri ohjelmointi (values/formulas first, then style) ---
$ws.Range("A11").Value = 46035
$ws.Range("B11").Value = 0.770833333333333
$ws.Range("C11").Value = 0.847222222222222
$ws.Range("D11").Formula = "=C11-B11"
$ws.Range("E11").Value = "Bäkkäri ohjelmointi"

$excel.CutCopyMode = $false

# --- Apply formatting to new/changed cells (after values/formulas are set) ---

# Date cells: A8, A9 (restyle) and A10, A11 (new) match A4:A7
$ws.Range("A4").Copy()
$ws.Range("A8:A11").PasteSpecial($fmt)

# Start/End time cells: B10:C11 match B4:C4
$ws.Range("B4:C4").Copy()
$ws.Range("B10:C11").PasteSpecial($fmt)

# Duration cells: D10:D11 match D4
$ws.Range("D4").Copy()
$ws.Range("D10:D11").PasteSpecial($fmt)

# Desc cells: E8, E9 (restyle) and E10 (new) match E4's format
# (General number format / default alignment - visually identical to the
# cells' previous default formatting, so this is a no-op where the source
# and destination formats already render the same way)
$ws.Range("E4").Copy()
$ws.Range("E8").PasteSpecial($fmt)
$ws.Range("E4").Copy()
$ws.Range("E9").PasteSpecial($fmt)
$ws.Range("E4").Copy()
$ws.Range("E10").PasteSpecial($fmt)
# (E11 format was already applied above, copied from the original E8 style)

$excel.CutCopyMode = $false

# Recalculate so cached formula values (e.g. SUM(D:D) in B1) are up to date
$excel.Calculate()

# Selection matches author's final cursor position
$ws.Range("C11").Select() | Out-Null
